$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 7 data (mirrors style/format of existing rows)
$ws.Range("A7").Value = 44226
$ws.Range("A7").NumberFormat = $ws.Range("A6").NumberFormat

$ws.Range("B7").Value = "PvZ "
$ws.Range("C7").Value = "Enh"
$ws.Range("E7").Value = "HIGH IMPORTANCE: Put a message that people make sure processess has been listed correctly before the OEdb is generated."
$ws.Range("D7").Value = "Prioritised"

$ws.Rows.Item(7).RowHeight = 43.5
$ws.Columns.Item(5).ColumnWidth = 42.5

$ws.Range("D8").Select()
